$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5042819496927109
$ws.Cells.Item(2, 3).Value = 0.2208144300982777
$ws.Cells.Item(2, 5).Value = 0.1266197142702836
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.5236478538078444
$ws.Cells.Item(2, 8).Value = 0.6775230762775024
$ws.Cells.Item(2, 11).Value = 0.2531345300494365
$ws.Cells.Item(2, 12).Value = 0.18473605578302
$ws.Cells.Item(2, 14).Value = 1.494652620132626
$ws.Cells.Item(2, 15).Value = 2.370649466758522

$ws.Cells.Item(3, 2).Value = 0.4653848426645766
$ws.Cells.Item(3, 3).Value = 0.2222145604939225
$ws.Cells.Item(3, 5).Value = 0.1260223839050916
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.5279957429694946
$ws.Cells.Item(3, 8).Value = 0.6828266683004927
$ws.Cells.Item(3, 11).Value = 0.2217909150135142
$ws.Cells.Item(3, 12).Value = 0.1776076767126114
$ws.Cells.Item(3, 14).Value = 1.508041058822563
$ws.Cells.Item(3, 15).Value = 2.390753126650097

$ws.Cells.Item(4, 2).Value = 0.4416086195346622
$ws.Cells.Item(4, 3).Value = 0.2231267684195046
$ws.Cells.Item(4, 5).Value = 0.1257220458618171
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.530988532030662
$ws.Cells.Item(4, 8).Value = 0.686341780868247
$ws.Cells.Item(4, 11).Value = 0.2025004094686977
$ws.Cells.Item(4, 12).Value = 0.1733237322842882
$ws.Cells.Item(4, 14).Value = 1.5167465026808
$ws.Cells.Item(4, 15).Value = 2.404316018096949

$ws.Cells.Item(5, 2).Value = 0.4319471982964558
$ws.Cells.Item(5, 3).Value = 0.2235117351988301
$ws.Cells.Item(5, 5).Value = 0.1256163904469432
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.5322893599235385
$ws.Cells.Item(5, 8).Value = 0.6878393249940444
$ws.Cells.Item(5, 11).Value = 0.1946284931477464
$ws.Cells.Item(5, 12).Value = 0.171601464190104
$ws.Cells.Item(5, 14).Value = 1.520416053392044
$ws.Cells.Item(5, 15).Value = 2.410149647956345

$ws.Cells.Item(6, 2).Value = 0.4303446142882308
$ws.Cells.Item(6, 3).Value = 0.2235764587877309
$ws.Cells.Item(6, 5).Value = 0.1255998581133895
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.532510267607023
$ws.Cells.Item(6, 8).Value = 0.6880919250974031
$ws.Cells.Item(6, 11).Value = 0.1933207262069914
$ws.Cells.Item(6, 12).Value = 0.1713169039290108
$ws.Cells.Item(6, 14).Value = 1.521032751157929
$ws.Cells.Item(6, 15).Value = 2.411136840298127

$ws.Cells.Item(7, 2).Value = 0.4414782095700218
$ws.Cells.Item(7, 3).Value = 0.2231319065935278
$ws.Cells.Item(7, 5).Value = 0.1257205531553822
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.5310057465322302
$ws.Cells.Item(7, 8).Value = 0.6863617135775399
$ws.Cells.Item(7, 11).Value = 0.2023942893879536
$ws.Cells.Item(7, 12).Value = 0.1733004099814508
$ws.Cells.Item(7, 14).Value = 1.516795497393435
$ws.Cells.Item(7, 15).Value = 2.404393450791872

$ws.Cells.Item(8, 2).Value = 0.4908484996780942
$ws.Cells.Item(8, 3).Value = 0.2212863106680967
$ws.Cells.Item(8, 5).Value = 0.12639998668703
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.5250799196658562
$ws.Cells.Item(8, 8).Value = 0.6792980999119465
$ws.Cells.Item(8, 11).Value = 0.2423369871155785
$ws.Cells.Item(8, 12).Value = 0.1822589716550169
$ws.Cells.Item(8, 14).Value = 1.499168365497436
$ws.Cells.Item(8, 15).Value = 2.377328214929719

$ws.Cells.Item(9, 2).Value = 0.5884830416406146
$ws.Cells.Item(9, 3).Value = 0.2180825892184544
$ws.Cells.Item(9, 5).Value = 0.1282582395984733
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 0.5160247508530205
$ws.Cells.Item(9, 8).Value = 0.6674963284218833
$ws.Cells.Item(9, 11).Value = 0.3202843136758986
$ws.Cells.Item(9, 12).Value = 0.2005600431948409
$ws.Cells.Item(9, 14).Value = 1.468444982246787
$ws.Cells.Item(9, 15).Value = 2.333924146883263

$ws.Cells.Item(10, 2).Value = 0.660684961709137
$ws.Cells.Item(10, 3).Value = 0.2159803544268719
$ws.Cells.Item(10, 5).Value = 0.1299428469724013
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 0.5109377686519991
$ws.Cells.Item(10, 8).Value = 0.6600717386517232
$ws.Cells.Item(10, 11).Value = 0.3772989872168182
$ws.Cells.Item(10, 12).Value = 0.2144493223464536
$ws.Cells.Item(10, 14).Value = 1.448209592742352
$ws.Cells.Item(10, 15).Value = 2.307927773902719

$ws.Cells.Item(11, 2).Value = 0.6936275992262608
$ws.Cells.Item(11, 3).Value = 0.2150782408726961
$ws.Cells.Item(11, 5).Value = 0.1307783164475005
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 0.5089640175233257
$ws.Cells.Item(11, 8).Value = 0.6569639552394975
$ws.Cells.Item(11, 11).Value = 0.4031772148407242
$ws.Cells.Item(11, 12).Value = 0.2208635360514393
$ws.Cells.Item(11, 14).Value = 1.439510224517832
$ws.Cells.Item(11, 15).Value = 2.297380323059002

$ws.Cells.Item(12, 2).Value = 0.7061154958975635
$ws.Cells.Item(12, 3).Value = 0.2147444004472234
$ws.Cells.Item(12, 5).Value = 0.1311046007334689
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 0.5082655826328519
$ws.Cells.Item(12, 8).Value = 0.6558258458589137
$ws.Cells.Item(12, 11).Value = 0.4129678007065252
$ws.Cells.Item(12, 12).Value = 0.2233061302422783
$ws.Cells.Item(12, 14).Value = 1.43628866115106
$ws.Cells.Item(12, 15).Value = 2.293570072974376

$ws.Cells.Item(13, 2).Value = 0.7034254280827668
$ws.Cells.Item(13, 3).Value = 0.2148159537454646
$ws.Cells.Item(13, 5).Value = 0.1310338893168499
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 0.5084138239572482
$ws.Cells.Item(13, 8).Value = 0.6560692357342646
$ws.Cells.Item(13, 11).Value = 0.4108596297810152
$ws.Cells.Item(13, 12).Value = 0.2227794678108523
$ws.Cells.Item(13, 14).Value = 1.4369792497329
$ws.Cells.Item(13, 15).Value = 2.294382502087004

$ws.Cells.Item(14, 2).Value = 0.6946547253577364
$ws.Cells.Item(14, 3).Value = 0.2150506200261688
$ws.Cells.Item(14, 5).Value = 0.1308049616172866
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 0.5089055749631868
$ws.Cells.Item(14, 8).Value = 0.6568695461121266
$ws.Cells.Item(14, 11).Value = 0.4039828747409615
$ws.Cells.Item(14, 12).Value = 0.2210642164619685
$ws.Cells.Item(14, 14).Value = 1.439243727774169
$ws.Cells.Item(14, 15).Value = 2.297063167030672

$ws.Cells.Item(15, 2).Value = 0.6892841164073786
$ws.Cells.Item(15, 3).Value = 0.2151953711254748
$ws.Cells.Item(15, 5).Value = 0.1306660264549713
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 0.5092131667389381
$ws.Cells.Item(15, 8).Value = 0.6573648032818582
$ws.Cells.Item(15, 11).Value = 0.399769482241993
$ws.Cells.Item(15, 12).Value = 0.2200153521988994
$ws.Cells.Item(15, 14).Value = 1.440640252414752
$ws.Cells.Item(15, 15).Value = 2.298729093296373

$ws.Cells.Item(16, 2).Value = 0.658533978631425
$ws.Cells.Item(16, 3).Value = 0.2160403983271877
$ws.Cells.Item(16, 5).Value = 0.1298896355561077
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 0.5110736100093831
$ws.Cells.Item(16, 8).Value = 0.6602802630800682
$ws.Cells.Item(16, 11).Value = 0.3756065657789804
$ws.Cells.Item(16, 12).Value = 0.2140320580244719
$ws.Cells.Item(16, 14).Value = 1.448788291966572
$ws.Cells.Item(16, 15).Value = 2.30864280025304

$ws.Cells.Item(17, 2).Value = 0.6396941971044328
$ws.Cells.Item(17, 3).Value = 0.2165726608391942
$ws.Cells.Item(17, 5).Value = 0.1294310294532188
$ws.Cells.Item(17, 6).Value = 0.6400460337125793
$ws.Cells.Item(17, 7).Value = 0.5123021309401423
$ws.Cells.Item(17, 8).Value = 0.6621378500119235
$ws.Cells.Item(17, 11).Value = 0.3607681140563557
$ws.Cells.Item(17, 12).Value = 0.2103859799249506
$ws.Cells.Item(17, 14).Value = 1.453916374146083
$ws.Cells.Item(17, 15).Value = 2.315051977939831

$ws.Cells.Item(18, 2).Value = 0.6288672897312324
$ws.Cells.Item(18, 3).Value = 0.2168839074809839
$ws.Cells.Item(18, 5).Value = 0.1291737603646261
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 0.5130407759291558
$ws.Cells.Item(18, 8).Value = 0.6632316727742165
$ws.Cells.Item(18, 11).Value = 0.3522280021059885
$ws.Cells.Item(18, 12).Value = 0.2082978871221428
$ws.Cells.Item(18, 14).Value = 1.456913522836668
$ws.Cells.Item(18, 15).Value = 2.318858699078362

$ws.Cells.Item(19, 2).Value = 0.6252030957989234
$ws.Cells.Item(19, 3).Value = 0.2169901674805033
$ws.Cells.Item(19, 5).Value = 0.1290877723459154
$ws.Cells.Item(19, 6).Value = 0.619163680173358
$ws.Cells.Item(19, 7).Value = 0.5132963687412797
$ws.Cells.Item(19, 8).Value = 0.6636063840386655
$ws.Cells.Item(19, 11).Value = 0.349335552895667
$ws.Cells.Item(19, 12).Value = 0.2075924502092192
$ws.Cells.Item(19, 14).Value = 1.457936484350668
$ws.Cells.Item(19, 15).Value = 2.320168256391739

$ws.Cells.Item(20, 2).Value = 0.6416987738654711
$ws.Cells.Item(20, 3).Value = 0.216515472586611
$ws.Cells.Item(20, 5).Value = 0.1294791754224462
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 0.512168037274094
$ws.Cells.Item(20, 8).Value = 0.6619374794633259
$ws.Cells.Item(20, 11).Value = 0.3623482585352349
$ws.Cells.Item(20, 12).Value = 0.2107731769622347
$ws.Cells.Item(20, 14).Value = 1.453365553708807
$ws.Cells.Item(20, 15).Value = 2.314357256128091

$ws.Cells.Item(21, 2).Value = 0.6972305405673751
$ws.Cells.Item(21, 3).Value = 0.2149814821568796
$ws.Cells.Item(21, 5).Value = 0.1308719345231033
$ws.Cells.Item(21, 6).Value = 0.7228739723491628
$ws.Cells.Item(21, 7).Value = 0.5087598060964496
$ws.Cells.Item(21, 8).Value = 0.6566334244726377
$ws.Cells.Item(21, 11).Value = 0.4060029916663268
$ws.Cells.Item(21, 12).Value = 0.2215676573134715
$ws.Cells.Item(21, 14).Value = 1.438576622943216
$ws.Cells.Item(21, 15).Value = 2.296270801346125

$ws.Cells.Item(22, 2).Value = 0.7336005298055568
$ws.Cells.Item(22, 3).Value = 0.214024211552811
$ws.Cells.Item(22, 5).Value = 0.1318399323244179
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 0.5068178521533824
$ws.Cells.Item(22, 8).Value = 0.6533927094787586
$ws.Cells.Item(22, 11).Value = 0.4344815151584385
$ws.Cells.Item(22, 12).Value = 0.2287021005824386
$ws.Cells.Item(22, 14).Value = 1.429334929365865
$ws.Cells.Item(22, 15).Value = 2.28552182288594

$ws.Cells.Item(23, 2).Value = 0.7141824231125042
$ws.Cells.Item(23, 3).Value = 0.2145309895753442
$ws.Cells.Item(23, 5).Value = 0.1313180199620376
$ws.Cells.Item(23, 6).Value = 0.7472568307830727
$ws.Cells.Item(23, 7).Value = 0.5078281708403409
$ws.Cells.Item(23, 8).Value = 0.655101693838489
$ws.Cells.Item(23, 11).Value = 0.4192869785402706
$ws.Cells.Item(23, 12).Value = 0.2248870683566224
$ws.Cells.Item(23, 14).Value = 1.434228635543995
$ws.Cells.Item(23, 15).Value = 2.291160705552272

$ws.Cells.Item(24, 2).Value = 0.6407924911503642
$ws.Cells.Item(24, 3).Value = 0.2165413110654271
$ws.Cells.Item(24, 5).Value = 0.129457388724628
$ws.Cells.Item(24, 6).Value = 0.6416283278902171
$ws.Cells.Item(24, 7).Value = 0.5122285602503638
$ws.Cells.Item(24, 8).Value = 0.6620279863981438
$ws.Cells.Item(24, 11).Value = 0.3616339041037691
$ws.Cells.Item(24, 12).Value = 0.210598099986413
$ws.Cells.Item(24, 14).Value = 1.453614427162545
$ws.Cells.Item(24, 15).Value = 2.31467095984182

$ws.Cells.Item(25, 2).Value = 0.5619856520685005
$ws.Cells.Item(25, 3).Value = 0.218904981958115
$ws.Cells.Item(25, 5).Value = 0.1276993613218522
$ws.Cells.Item(25, 6).Value = 0.5279251897347166
$ws.Cells.Item(25, 7).Value = 0.5181995726129216
$ws.Cells.Item(25, 8).Value = 0.6704699254214717
$ws.Cells.Item(25, 11).Value = 0.2992404166050733
$ws.Cells.Item(25, 12).Value = 0.1955309763958439
$ws.Cells.Item(25, 14).Value = 1.476345635886197
$ws.Cells.Item(25, 15).Value = 2.344631080348435

Write-Host "Updated 240 cells"